$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated crypto price/volume data scraped on
# Wed May 22 08:48:06 UTC 2024.
#
# Price-column (D) values are plain text (e.g. thousands are dot-
# separated, "70.043.84"), so any replacement value that Excel would
# otherwise auto-parse as a number (a single decimal point, no extra
# separators) needs its cell pre-formatted as Text ("@") to keep it a
# string literal instead of silently becoming a numeric value.

$ws.Range("D2").Value = '70.043.84'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '3.767.81'
$ws.Range("E3").Value = '  +3.17%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '624.41'
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.94'
$ws.Range("E6").Value = '  -0.78%  '
$ws.Range("D7").Value = '3.765.79'
$ws.Range("E7").Value = '  +3.20%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("E10").Value = '  +3.09%  '
$ws.Range("E11").Value = '  -5.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.489'
$ws.Range("E12").Value = '  -2.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.21'
$ws.Range("E13").Value = '  +1.98%  '
$ws.Range("D15").Value = '4.387.87'
$ws.Range("E15").Value = '  +2.88%  '
$ws.Range("D16").Value = '3.766.48'
$ws.Range("E16").Value = '  +2.91%  '
$ws.Range("D17").Value = '70.135.33'
$ws.Range("E17").Value = '  -1.08%  '
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.66'
$ws.Range("E19").Value = '  +1.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.80'
$ws.Range("E20").Value = '  -0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '507.42'
$ws.Range("E21").Value = '  -2.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.48'
$ws.Range("E22").Value = '  +2.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.729'
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.52'
$ws.Range("E24").Value = '  -1.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.25'
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("E26").Value = '  -1.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.18'
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("E28").Value = '  +26.07%  '
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("E30").Value = '  -1.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.97'
$ws.Range("E31").Value = '  +2.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.90'
$ws.Range("E32").Value = '  -3.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.54'
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("E34").Value = '  +0.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  +4.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.23'
$ws.Range("E37").Value = '  +1.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.336'
$ws.Range("E38").Value = '  -4.02%  '
$ws.Range("E39").Value = '  +0.97%  '
$ws.Range("E40").Value = '  -3.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.34'
$ws.Range("E41").Value = '  -2.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.23'
$ws.Range("E42").Value = '  -1.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '423.97'
$ws.Range("E43").Value = '  -0.70%  '
$ws.Range("B44").Value = 'Cosmos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.75'
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.86'
$ws.Range("E45").Value = '  +3.05%  '
$ws.Range("D46").Value = '3.012.57'
$ws.Range("E46").Value = '  -3.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0367'
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.43'
$ws.Range("E48").Value = '  -3.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.51'
$ws.Range("E49").Value = '  -1.42%  '
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.54'
$ws.Range("E51").Value = '  +2.68%  '
